$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 203, shifting existing rows 203:237 down to 204:238
$ws.Rows.Item(203).Insert()

# Populate the newly inserted row 203 with the new record's data
$ws.Cells.Item(203, 1).Value = 3
$ws.Cells.Item(203, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(203, 3).Value = "Coquimbo"
$ws.Cells.Item(203, 4).Value = 44522
$ws.Cells.Item(203, 5).Value = 5
$ws.Cells.Item(203, 6).Value = 100112043
$ws.Cells.Item(203, 7).Value = "Pepino ensalada"
$ws.Cells.Item(203, 8).Value = "Sin especificar"
$ws.Cells.Item(203, 9).Value = "Primera"
$ws.Cells.Item(203, 10).Value = 170
$ws.Cells.Item(203, 11).Value = 7000
$ws.Cells.Item(203, 12).Value = 7500
$ws.Cells.Item(203, 13).Value = 7235
$ws.Cells.Item(203, 14).Value = "`$/caja 70 unidades"
$ws.Cells.Item(203, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(203, 16).Value = 103
$ws.Cells.Item(203, 17).Value = 70
$ws.Cells.Item(203, 18).Value = "Hortaliza"
